$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Sdate column (D2:D10) to the new date string for every event row.
$ws.Range("D2:D10").Value = "2023-11-01T11:25:00"

# Update the selection to match the new state (D2:D10 selected, active cell D2)
$ws.Range("D2:D10").Select()
